$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.011.17"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -4.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.222.90"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  -5.69%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.33"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.16"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -8.84%  "

$ws.Range("E7").Value = "  -7.38%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -8.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.44"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -11.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.30"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -2.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -9.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.72"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -8.51%  "

$ws.Range("E14").Value = "  -3.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.868"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  -11.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.563.53"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  -5.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.02"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -8.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.217.51"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  -5.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.871.25"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -5.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.59"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +4.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0965"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -9.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -11.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.21"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  -10.83%  "

$ws.Range("E24").Value = "  -8.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.16"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -9.00%  "

$ws.Range("E26").Value = "  -8.88%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  -8.95%  "

$ws.Range("E29").Value = "  -6.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.27"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -14.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0882"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -8.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.48"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -8.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.74"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  -10.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.36"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -7.29%  "

$ws.Range("E35").Value = "  -5.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.33"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +11.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.04"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +18.22%  "

$ws.Range("E38").Value = "  -6.17%  "

$ws.Range("E39").Value = "  -7.21%  "

$ws.Range("E40").Value = "  -11.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -6.69%  "

$ws.Range("E42").Value = "  -8.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.879.04"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +11.79%  "

$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.10"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -5.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.79"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -11.61%  "

$ws.Range("E47").Value = "  -10.93%  "

$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "78.31"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -4.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.51"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -12.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.62"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  -5.82%  "
